$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 80
$prev = $row - 1

$ws.Cells.Item($row, 1).Value = 79
$ws.Cells.Item($row, 2).Value = "denmark"
$ws.Cells.Item($row, 3).Value = "superliga"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45234.6875
$ws.Cells.Item($row, 6).Value = "Lyngby"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Odense"
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 2.06
$ws.Cells.Item($row, 11).Value = "29/10/2023 16:12"
$ws.Cells.Item($row, 12).Value = 2.49
$ws.Cells.Item($row, 13).Value = "04/11/2023 16:20"
$ws.Cells.Item($row, 14).Value = 3.66
$ws.Cells.Item($row, 15).Value = "29/10/2023 16:12"
$ws.Cells.Item($row, 16).Value = 3.53
$ws.Cells.Item($row, 17).Value = "04/11/2023 16:20"
$ws.Cells.Item($row, 18).Value = 3.5
$ws.Cells.Item($row, 19).Value = "29/10/2023 16:12"
$ws.Cells.Item($row, 20).Value = 2.88
$ws.Cells.Item($row, 21).Value = "04/11/2023 16:20"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/denmark/superliga/lyngby-odense/z5BdBSjA/"

# Column A (Indice) carries the bold/border/centered style used throughout the table.
$ws.Cells.Item($prev, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)

# Column E (data_partida) carries the custom date-time number format.
$ws.Cells.Item($prev, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)

$excel.CutCopyMode = 0
